$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2018-12-31 00:00:00"
$ws.Range("O2").Value = 45459286.49
$ws.Range("P2").Value = 426025451.47
$ws.Range("Q2").Value = 366187885.59
$ws.Range("R2").Value = 103.0933612876
$ws.Range("S2").Value = 297075662.62
$ws.Range("T2").Value = 297075662.62
$ws.Range("U2").Value = 159.8546501352
$ws.Range("V2").Value = 22308039.16
$ws.Range("W2").Value = 20388931.39
$ws.Range("X2").Value = -1969925.65
$ws.Range("Y2").Value = 50938396.88
$ws.Range("Z2").Value = 50683187.23
$ws.Range("AA2").Value = 5223900.74
$ws.Range("AG2").Value = 1804463.77
$ws.Range("AP2").Value = 197.2180718781
$ws.Range("AQ2").Value = 194.818148127354
$ws.Range("AR2").Value = 203.023809830422
$ws.Range("AS2").Value = 41137086.49
$ws.Range("AT2").Value = 188.46092818582
